$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.480.88'
$ws.Range('E2').Value = '  +1.81%  '
$ws.Range('D3').Value = '1.866.60'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '311.85'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('E6').Value = '  -0.06%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4779'
$ws.Range('D7').NumberFormat = 'General'
$ws.Range('E7').Value = '  +0.19%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3767'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('E8').Value = '  +2.48%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.07337'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('E9').Value = '  +1.52%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.9363'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('E10').Value = '  +0.64%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '20.70'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +4.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07843'
$ws.Range('D12').NumberFormat = 'General'
$ws.Range('E12').Value = '  +1.49%  '
$ws.Range('D13').Value = '1.896.92'
$ws.Range('E13').Value = '  +1.52%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.442'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +2.18%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.555'
$ws.Range('D15').NumberFormat = 'General'
$ws.Range('E15').Value = '  +2.08%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '90.40'
$ws.Range('D16').NumberFormat = 'General'
$ws.Range('E16').Value = '  +1.72%  '
$ws.Range('E17').Value = '  -0.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.000008892'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +2.86%  '
$ws.Range('E19').Value = '  -0.08%  '
$ws.Range('D20').Value = '27.541.63'
$ws.Range('E20').Value = '  +1.91%  '
$ws.Range('E21').Value = '  +1.22%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.122'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +1.17%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '10.69'
$ws.Range('D23').NumberFormat = 'General'
$ws.Range('E23').Value = '  +0.27%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.940'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '154.49'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +1.14%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '18.47'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.023'
$ws.Range('D27').NumberFormat = 'General'
$ws.Range('E27').Value = '  +1.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '115.51'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('E28').Value = '  +0.93%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '4.972'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  -0.29%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.08902'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +0.16%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.329'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.02%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.217'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +3.78%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.7605'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +2.05%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.610'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +2.36%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.758'
$ws.Range('D35').NumberFormat = 'General'
$ws.Range('E35').Value = '  +0.96%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.125'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +1.19%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02038'
$ws.Range('D37').NumberFormat = 'General'
$ws.Range('E37').Value = '  +4.16%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.995'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.05262'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  -0.32%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.5316'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  +2.22%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.084'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +0.96%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.518'
$ws.Range('D42').NumberFormat = 'General'
$ws.Range('E42').Value = '  +3.76%  '
$ws.Range('E43').Value = '  +0.92%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.67'
$ws.Range('D44').NumberFormat = 'General'
$ws.Range('E44').Value = '  +0.20%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.4805'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +1.30%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.012'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('E46').Value = '  -0.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '102.94'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +1.31%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.653'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +2.84%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '67.38'
$ws.Range('D49').NumberFormat = 'General'
$ws.Range('E49').Value = '  +2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06079'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9175'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +3.27%  '
